$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("A4").Value = 112223199
$ws.Range("B4").Value = 57584
$ws.Range("E4").Value = 208245
$ws.Range("F4").Value = "Vanlig padda"
$ws.Range("G4").Value = "Bufo bufo"
$ws.Range("I4").NumberFormat = "@"
$ws.Range("I4").Value = "1"
$ws.Range("M4").Value = "spel/sång"
$ws.Range("Q4").Value = 636422
$ws.Range("R4").Value = 6523909
$ws.Range("AC4").Value = "Observerad med ljud."

# Row 5
$ws.Range("A5").Value = 112223196
$ws.Range("B5").Value = 57620
$ws.Range("E5").Value = 208242
$ws.Range("F5").Value = "Mindre vattensalamander"
$ws.Range("G5").Value = "Lissotriton vulgaris"
$ws.Range("H5").Value = "(Linnaeus, 1758)"
$ws.Range("J5").Value = "ex."
$ws.Range("K5").Value = "adult"
$ws.Range("L5").Value = "hane"
$ws.Range("M5").Value = "i vatten/simmande"
$ws.Range("Q5").Value = 636422
$ws.Range("R5").Value = 6523909
$ws.Range("Z5").Value = "20:30"
$ws.Range("AB5").Value = "20:30"
$ws.Range("AC5").Value = "Hane i lekdräkt."

# Row 6
$ws.Range("A6").Value = 112223184
$ws.Range("B6").Value = 57610
$ws.Range("E6").Value = 208250
$ws.Range("F6").Value = "Åkergroda"
$ws.Range("G6").Value = "Rana arvalis"
$ws.Range("H6").Value = "Nilsson, 1842"
$ws.Range("J6").Value = "äggklumpar"
$ws.Range("K6").Value = "ägg"
$ws.Range("L6").Value = ""
$ws.Range("M6").Value = ""
$ws.Range("AC6").Value = "Romklump. Troligen åkergroda som förekommer i närliggande dammar."

# Row 7
$ws.Range("A7").Value = 112223188
$ws.Range("B7").Value = 57620
$ws.Range("Q7").Value = 636399
$ws.Range("R7").Value = 6523963
$ws.Range("Z7").Value = "20:20"
$ws.Range("AB7").Value = "20:20"
$ws.Range("AC7").Value = ""

# Row 8
$ws.Range("A8").Value = 112223201
$ws.Range("B8").Value = 57584
$ws.Range("E8").Value = 208245
$ws.Range("F8").Value = "Vanlig padda"
$ws.Range("G8").Value = "Bufo bufo"
$ws.Range("H8").Value = "(Linnaeus, 1758)"
$ws.Range("I8").NumberFormat = "@"
$ws.Range("I8").Value = "2"
$ws.Range("AC8").Value = "Två hanar observerade med ljud och visuellt."

# Row 9
$ws.Range("A9").Value = 112223193
$ws.Range("B9").Value = 57620
$ws.Range("E9").Value = 208242
$ws.Range("F9").Value = "Mindre vattensalamander"
$ws.Range("G9").Value = "Lissotriton vulgaris"
$ws.Range("I9").NumberFormat = "@"
$ws.Range("I9").Value = "5"
$ws.Range("M9").Value = "i vatten/simmande"
$ws.Range("Q9").Value = 636408
$ws.Range("R9").Value = 6524025
$ws.Range("Z9").Value = "20:30"
$ws.Range("AB9").Value = "20:30"
$ws.Range("AC9").Value = "Hanar i lekdräkt."

# Row 10
$ws.Range("A10").Value = 112223203
$ws.Range("B10").Value = 57610
$ws.Range("E10").Value = 208250
$ws.Range("F10").Value = "Åkergroda"
$ws.Range("G10").Value = "Rana arvalis"
$ws.Range("H10").Value = "Nilsson, 1842"
$ws.Range("Q10").Value = 636549
$ws.Range("R10").Value = 6523814
$ws.Range("Z10").Value = "21:15"
$ws.Range("AB10").Value = "21:15"
